# Updated symbol list on Mon Feb 13 14:33:06 UTC 2023 with GitHub Actions
# Refresh Price (column D) and Volume(1h) (column E) figures for the crypto
# symbol table on the active sheet. Source cells are plain text (e.g. "289.75",
# "-6.58%"), not numbers/percentages, so each write forces a Text number format
# before assigning the value (otherwise Excel auto-coerces the literal into a
# numeric/percentage value) and then restores the default "Normal" style so no
# extra cell formatting is introduced.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($cellRef, $value) {
    $rng = $ws.Range($cellRef)
    $rng.NumberFormat = "@"
    $rng.Value = $value
    $rng.Style = "Normal"
}

Set-TextValue "D2" "289.75"
Set-TextValue "E2" "-6.58%"
Set-TextValue "D3" "39.70"
Set-TextValue "E3" "-3.36%"
Set-TextValue "E4" "-3.70%"
Set-TextValue "D5" "0.07349"
Set-TextValue "E5" "-4.07%"
Set-TextValue "E6" "-0.39%"
Set-TextValue "D7" "1.553"
Set-TextValue "E7" "-11.23%"
Set-TextValue "D8" "0.9156"
Set-TextValue "E8" "-0.57%"
Set-TextValue "D9" "0.1188"
Set-TextValue "E9" "-6.25%"
Set-TextValue "D10" "0.1732"
Set-TextValue "E10" "-4.90%"
Set-TextValue "D11" "0.08677"
Set-TextValue "E11" "-4.98%"
Set-TextValue "D12" "0.04179"
Set-TextValue "E12" "0.48%"
Set-TextValue "D13" "0.1050"
Set-TextValue "E13" "-0.10%"
Set-TextValue "D14" "0.001279"
Set-TextValue "E14" "-0.44%"
Set-TextValue "D15" "0.005850"
Set-TextValue "E15" "-0.76%"
Set-TextValue "D16" "3.381"
Set-TextValue "E16" "0.81%"
Set-TextValue "D18" "0.3296"
Set-TextValue "E18" "-0.72%"
Set-TextValue "D19" "7.547"
Set-TextValue "E19" "1.85%"
Set-TextValue "D20" "0.1353"
Set-TextValue "E20" "-0.16%"
Set-TextValue "D21" "0.2735"
Set-TextValue "E21" "0.44%"
Set-TextValue "D22" "0.03842"
Set-TextValue "E22" "-4.11%"
Set-TextValue "D23" "0.001272"
Set-TextValue "E23" "0.27%"
Set-TextValue "D24" "0.003692"
Set-TextValue "E24" "-9.84%"
Set-TextValue "D25" "0.0001284"
Set-TextValue "E25" "0.89%"
Set-TextValue "D26" "0.0003730"
Set-TextValue "D38" "0.02316"
Set-TextValue "E38" "-7.81%"
Set-TextValue "D39" "0.05005"
Set-TextValue "E39" "-5.70%"
Set-TextValue "D40" "0.007703"
Set-TextValue "E40" "-1.95%"
Set-TextValue "E41" "148.63%"
Set-TextValue "D42" "0.1270"
Set-TextValue "E42" "-3.04%"
Set-TextValue "D43" "0.007407"
Set-TextValue "E43" "11.14%"
Set-TextValue "D44" "0.007697"
Set-TextValue "E44" "-5.24%"
Set-TextValue "D45" "0.3151"
Set-TextValue "E45" "2.30%"
Set-TextValue "D46" "0.00006510"
Set-TextValue "E46" "-4.25%"
Set-TextValue "E47" "0.09%"
Set-TextValue "E48" "12.37%"
Set-TextValue "D49" "0.004209"
Set-TextValue "E49" "35.62%"
Set-TextValue "E50" "0.09%"
Set-TextValue "E51" "0.09%"
